$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("A1").Value = "First_Name"
$ws.Range("B1").Value = "Last_Name"
$ws.Range("C1").Value = "Postal_Code"
$ws.Range("C2").Select()
